$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row of test-case data (row 12)
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "e -f encrypted.txt -k abc"

$ws.Range("D12").Value = "sihT@si@elpmas@elif"
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:sihT@si@elpmas@elif")

$ws.Range("E12").Value = "that is encrypted file"
$ws.Range("F12").Value = "FAIL"

$ws.Range("F13").Select()
